$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2097167003751395
$ws.Range("D2").Value = 0.1683836165379375
$ws.Range("E2").Value = 0.1572518775516691
$ws.Range("F2").Value = 1.585520969950849
$ws.Range("G2").Value = 0.9885403027152648
$ws.Range("H2").Value = 0.9588152100885736
$ws.Range("I2").Value = 1.107378677071821
$ws.Range("J2").Value = 0.1853140390250303
$ws.Range("L2").Value = 0.2123279231335857
$ws.Range("O2").Value = 3.943819789077736

$ws.Range("C3").Value = 0.2089922281128764
$ws.Range("D3").Value = 0.1682660490980226
$ws.Range("E3").Value = 0.1569064231002919
$ws.Range("F3").Value = 1.562029610908425
$ws.Range("G3").Value = 0.9632484458721677
$ws.Range("H3").Value = 0.9514223883856232
$ws.Range("I3").Value = 1.090128063910129
$ws.Range("J3").Value = 0.1846786843179231
$ws.Range("L3").Value = 0.2120131430796448
$ws.Range("O3").Value = 3.87433464924348

$ws.Range("C4").Value = 0.2086295553522319
$ws.Range("D4").Value = 0.1682508060830017
$ws.Range("E4").Value = 0.1567542179741217
$ws.Range("F4").Value = 1.5484108024786
$ws.Range("G4").Value = 0.9482747664937108
$ws.Range("H4").Value = 0.9473179711940389
$ws.Range("I4").Value = 1.080104027687696
$ws.Range("J4").Value = 0.1843640825306352
$ws.Range("L4").Value = 0.2118993815216257
$ws.Range("O4").Value = 3.833743338436676

$ws.Range("C5").Value = 0.2085024639902073
$ws.Range("D5").Value = 0.1682589471707132
$ws.Range("E5").Value = 0.1567072846196673
$ws.Range("F5").Value = 1.543063623931715
$ws.Range("G5").Value = 0.9423126072115622
$ws.Range("H5").Value = 0.9457548628083856
$ws.Range("I5").Value = 1.076162072721239
$ws.Range("J5").Value = 0.1842548947162967
$ws.Range("L5").Value = 0.2118730520793477
$ws.Range("O5").Value = 3.817723608565984

$ws.Range("C6").Value = 0.2084826120744765
$ws.Range("D6").Value = 0.1682611668560554
$ws.Range("E6").Value = 0.1567004036798032
$ws.Range("F6").Value = 1.542187970625051
$ws.Range("G6").Value = 0.9413310356306255
$ws.Range("H6").Value = 0.9455019260073669
$ws.Range("I6").Value = 1.075516149641331
$ws.Range("J6").Value = 0.1842379134400858
$ws.Range("L6").Value = 0.2118698908423582
$ws.Range("O6").Value = 3.815095055259434

$ws.Range("C7").Value = 0.2086277574820983
$ws.Range("D7").Value = 0.168250857717247
$ws.Range("E7").Value = 0.156753523872549
$ws.Range("F7").Value = 1.548337868022202
$ws.Range("G7").Value = 0.9481937929152764
$ws.Range("H7").Value = 0.9472964471265044
$ws.Range("I7").Value = 1.080050286271614
$ws.Range("J7").Value = 0.1843625329575076
$ws.Range("L7").Value = 0.211898945289235
$ws.Range("O7").Value = 3.833525178567442

$ws.Range("C8").Value = 0.209449875280157
$ws.Range("D8").Value = 0.1683312804561581
$ws.Range("E8").Value = 0.1571203477302028
$ws.Range("F8").Value = 1.577254170375852
$ws.Range("G8").Value = 0.9797043527199634
$ws.Range("H8").Value = 0.956176000301781
$ws.Range("I8").Value = 1.101312856878231
$ws.Range("J8").Value = 0.1850793129122792
$ws.Range("L8").Value = 0.2122029039555215
$ws.Range("O8").Value = 3.919431255942129

$ws.Range("C9").Value = 0.2117123634848923
$ws.Range("D9").Value = 0.1689394810343856
$ws.Range("E9").Value = 0.1583139770020487
$ws.Range("F9").Value = 1.640343245647259
$ws.Range("G9").Value = 1.045908725391172
$ws.Range("H9").Value = 0.9770346320178476
$ws.Range("I9").Value = 1.147513289427891
$ws.Range("J9").Value = 0.1870830904516865
$ws.Range("L9").Value = 0.213428588197921
$ws.Range("O9").Value = 4.104339501691641

$ws.Range("C10").Value = 0.2137693990293172
$ws.Range("D10").Value = 0.1696593889176725
$ws.Range("E10").Value = 0.1594789656588453
$ws.Range("F10").Value = 1.690590273277067
$ws.Range("G10").Value = 1.097250556706655
$ws.Range("H10").Value = 0.9944575747461215
$ws.Range("I10").Value = 1.184206116318748
$ws.Range("J10").Value = 0.1889190421313245
$ws.Range("L10").Value = 0.2147115437643023
$ws.Range("O10").Value = 4.250237925732165

$ws.Range("C11").Value = 0.2147905982924527
$ws.Range("D11").Value = 0.1700458703943895
$ws.Range("E11").Value = 0.1600712770171526
$ws.Range("F11").Value = 1.714296040369845
$ws.Range("G11").Value = 1.121196801144606
$ws.Range("H11").Value = 1.002838971503138
$ws.Range("I11").Value = 1.201496678672271
$ws.Range("J11").Value = 0.1898330956424701
$ws.Range("L11").Value = 0.2153779694204232
$ws.Range("O11").Value = 4.318797991820588

$ws.Range("C12").Value = 0.2151895514552677
$ws.Range("D12").Value = 0.1702006735885604
$ws.Range("E12").Value = 0.160304513033541
$ws.Range("F12").Value = 1.723394706966047
$ws.Range("G12").Value = 1.130349641506569
$ws.Range("H12").Value = 1.006078217826712
$ws.Range("I12").Value = 1.208130248121691
$ws.Range("J12").Value = 0.1901905446573693
$ws.Range("L12").Value = 0.2156422058599787
$ws.Range("O12").Value = 4.345074874831425

$ws.Range("C13").Value = 0.2151030856858682
$ws.Range("D13").Value = 0.1701669585946988
$ws.Range("E13").Value = 0.1602538842570524
$ws.Range("F13").Value = 1.721429731369284
$ws.Range("G13").Value = 1.128374636983381
$ws.Range("H13").Value = 1.005377681741805
$ws.Range("I13").Value = 1.206697766567913
$ws.Range("J13").Value = 0.190113058689299
$ws.Range("L13").Value = 0.2155847701470677
$ws.Range("O13").Value = 4.339401686917313

$ws.Range("C14").Value = 0.2148231751840797
$ws.Range("D14").Value = 0.1700584369570208
$ws.Range("E14").Value = 0.1600902864487139
$ws.Range("F14").Value = 1.715042152834997
$ws.Range("G14").Value = 1.121948109009395
$ws.Range("H14").Value = 1.003104156372245
$ws.Range("I14").Value = 1.202040703248656
$ws.Range("J14").Value = 0.1898622765233213
$ws.Range("L14").Value = 0.2153994704876467
$ws.Range("O14").Value = 4.320953501781673

$ws.Range("C15").Value = 0.2146533156918053
$ws.Range("D15").Value = 0.1699930639515088
$ws.Range("E15").Value = 0.1599912416342448
$ws.Range("F15").Value = 1.711145434843814
$ws.Range("G15").Value = 1.118022733549907
$ws.Range("H15").Value = 1.001720068928165
$ws.Range("I15").Value = 1.199199315425702
$ws.Range("J15").Value = 0.1897101382578157
$ws.Range("L15").Value = 0.215287514708109
$ws.Range("O15").Value = 4.309694427688271

$ws.Range("C16").Value = 0.2137043772181073
$ws.Range("D16").Value = 0.1696353157801553
$ws.Range("E16").Value = 0.1594415092120869
$ws.Range("F16").Value = 1.689058099189737
$ws.Range("G16").Value = 1.095697496723346
$ws.Range("H16").Value = 0.9939189882891242
$ws.Range("I16").Value = 1.183088178493207
$ws.Range("J16").Value = 0.1888608915408696
$ws.Range("L16").Value = 0.2146696547924236
$ws.Range("O16").Value = 4.245801412653179

$ws.Range("C17").Value = 0.2131440901763426
$ws.Range("D17").Value = 0.1694309353381556
$ws.Range("E17").Value = 0.1591202178067057
$ws.Range("F17").Value = 1.675725353392465
$ws.Range("G17").Value = 1.082152944683912
$ws.Range("H17").Value = 0.9892498833325476
$ws.Range("I17").Value = 1.173357816738346
$ws.Range("J17").Value = 0.1883600890957524
$ws.Range("L17").Value = 0.2143118030946383
$ws.Range("O17").Value = 4.207165875450755

$ws.Range("C18").Value = 0.2128298725899356
$ws.Range("D18").Value = 0.1693189355992288
$ws.Range("E18").Value = 0.1589412898942584
$ws.Range("F18").Value = 1.668136550989459
$ws.Range("G18").Value = 1.074418066012811
$ws.Range("H18").Value = 0.9866072293752097
$ws.Range("I18").Value = 1.16781753888327
$ws.Range("J18").Value = 0.1880794658428115
$ws.Range("L18").Value = 0.2141137719156418
$ws.Range("O18").Value = 4.18514989848876

$ws.Range("C19").Value = 0.2127248668035406
$ws.Range("D19").Value = 0.1692819693728183
$ws.Range("E19").Value = 0.1588817170073575
$ws.Range("F19").Value = 1.665580830279879
$ws.Range("G19").Value = 1.071808716799637
$ws.Range("H19").Value = 0.9857198426240075
$ws.Range("I19").Value = 1.165951380452938
$ws.Range("J19").Value = 0.1879857276462147
$ws.Range("L19").Value = 0.214048061832429
$ws.Range("O19").Value = 4.177731083675781

$ws.Range("C20").Value = 0.2132029013800292
$ws.Range("D20").Value = 0.1694521173479018
$ws.Range("E20").Value = 0.1591538124464051
$ws.Range("F20").Value = 1.67713638488658
$ws.Range("G20").Value = 1.083589032134142
$ws.Range("H20").Value = 0.9897424792854395
$ws.Range("I20").Value = 1.174387797112047
$ws.Range("J20").Value = 0.1884126320832493
$ws.Range("L20").Value = 0.2143490904263174
$ws.Range("O20").Value = 4.211257359249259

$ws.Range("C21").Value = 0.214905059614054
$ws.Range("D21").Value = 0.170090083287171
$ws.Range("E21").Value = 0.1601380966499448
$ws.Range("F21").Value = 1.716915034907331
$ws.Range("G21").Value = 1.123833431492841
$ws.Range("H21").Value = 1.003770172107039
$ws.Range("I21").Value = 1.203406262172649
$ws.Range("J21").Value = 0.1899356303855129
$ws.Range("L21").Value = 0.2154535754410603
$ws.Range("O21").Value = 4.326363642893227

$ws.Range("C22").Value = 0.2160888860254317
$ws.Range("D22").Value = 0.1705562753562262
$ws.Range("E22").Value = 0.1608334839038825
$ws.Range("F22").Value = 1.743622595914971
$ws.Range("G22").Value = 1.150630432039065
$ws.Range("H22").Value = 1.013319185784411
$ws.Range("I22").Value = 1.222872785554358
$ws.Range("J22").Value = 0.1909969494175598
$ws.Range("L22").Value = 0.2162446253436059
$ws.Range("O22").Value = 4.403426292070321

$ws.Range("C23").Value = 0.215450538467266
$ws.Range("D23").Value = 0.1703029645135388
$ws.Range("E23").Value = 0.1604575838574505
$ws.Range("F23").Value = 1.729303361729592
$ws.Range("G23").Value = 1.136283082525438
$ws.Range("H23").Value = 1.008187866786017
$ws.Range("I23").Value = 1.21243730136905
$ws.Range("J23").Value = 0.1904244769857328
$ws.Range("L23").Value = 0.2158161051852687
$ws.Range("O23").Value = 4.362128778469867

$ws.Range("C24").Value = 0.2131762882276007
$ws.Range("D24").Value = 0.1694425238239532
$ws.Range("E24").Value = 0.1591386062829407
$ws.Range("F24").Value = 1.676498219639825
$ws.Range("G24").Value = 1.082939614816098
$ws.Range("H24").Value = 0.989519646818195
$ws.Range("I24").Value = 1.173921975274965
$ws.Range("J24").Value = 0.1883888546747983
$ws.Range("L24").Value = 0.2143322088291981
$ws.Range("O24").Value = 4.209406988518708

$ws.Range("C25").Value = 0.2110308521812527
$ws.Range("D25").Value = 0.1687268810623976
$ws.Range("E25").Value = 0.1579404079817479
$ws.Range("F25").Value = 1.622592143835874
$ws.Range("G25").Value = 1.027525122434014
$ws.Range("H25").Value = 0.9710231873908981
$ws.Range("I25").Value = 1.134532267181697
$ws.Range("J25").Value = 0.1864770691823665
$ws.Range("L25").Value = 0.2130297470715661
$ws.Range("O25").Value = 4.052554133090553
